# ---------------------------------------------------------------------------
# Applies the "AquiferOpenStudyNotesBookIntros" resource-data update:
#   1. Heading2 "License Information" -> Normal-style bold
#      "Aquifer Open Study Notes (Book Intros)"
#   2. License paragraph body text rewritten (Tyndale -> Aquifer adaptation
#      notice); the two inline hyperlinks (Tyndale site, CC license) are
#      removed and folded into plain text.
#   3. "This PDF version is provided under the same license." replaced with
#      the multi-language adaptation sentence.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. "License Information" heading -------------------------------------
$headingPara = $d.Paragraphs.Item(4)
$headingRange = $headingPara.Range
$headingRange.MoveEnd(1, -1) | Out-Null
$headingRange.Text = "Aquifer Open Study Notes (Book Intros)"
$headingRange.Font.Bold = 1
$headingPara.Style = "Normal"

# --- 2. License body paragraph ---------------------------------------------
$licensePara = $d.Paragraphs.Item(5)

# 2a. Collapse the bold "Notas de Estudo..." run + the two plain runs that
#     follow it into a single, non-bold lead-in sentence.
$r = $licensePara.Range
$r.Find.ClearFormatting()
$r.Find.Execute(
    "Notas de Estudo - Introduções aos Livros (Tyndale) (Portuguese) is based on: ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This work is an adaptation of ", 2) | Out-Null
$licensePara.Range.Font.Bold = 0

# 2b. Remove the two hyperlinks in this paragraph (Tyndale House Publishers,
#     CC BY-SA 4.0 license) and replace the surrounding text in one pass.
$r2 = $licensePara.Range
$r2.Find.ClearFormatting()
$r2.Find.MatchWildcards = $true
$r2.Find.Execute(
    ", Tyndale House Publishers, 2019, which is licensed under a CC BY-SA 4.0 license.",
    $true, $false, $true, $false, $false, $true, 1, $false,
    " © 2023 Tyndale House Publishers, licensed under the CC BY-SA 4.0 license. The adaptation, Aquifer Open Study Notes, was created by Mission Mutual and is also licensed under CC BY-SA 4.0.",
    2) | Out-Null
$licensePara.Range.Font.Bold = 0

# --- 3. "This PDF version..." sentence -------------------------------------
$d.Content.Find.ClearFormatting()
$d.Content.Find.Execute(
    "This PDF version is provided under the same license.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This resource has been adapted into multiple languages, including English, Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文).",
    2) | Out-Null
